$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "GNG_TO-16511687345106356"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687365088973"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687365108652"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651168736572327"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687366475823"

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168734470461.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687344938858.csv"
$ws1.Range("B4").Value = "go_stims-1651168734495889.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687345096183.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511687361683564.csv"
$ws2.Range("B3").Value = "TB-16511687364288855.csv"
$ws2.Range("B4").Value = "OB-16511687356424897.csv"
$ws2.Range("B5").Value = "OB-16511687359896805.csv"
$ws2.Range("B6").Value = "ZB-match_9-16511687347324631.csv"
$ws2.Range("B7").Value = "ZB-match_1-1651168735249095.csv"
$ws2.Range("B8").Value = "OB-1651168735686449.csv"
$ws2.Range("B9").Value = "TB-16511687364879131.csv"
$ws2.Range("B10").Value = "ZB-match_2-165116873519122.csv"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687365238693.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168736512958.csv"
$ws4.Range("B4").Value = "MM_stims-1651168736554606.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687365248315.csv"
$ws4.Range("B6").Value = "MM_stims-1651168736570327.csv"
$ws4.Range("B7").Value = "ZM_stims-1651168736554606.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687365766115.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651168736616561.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511687366325836.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687366010346.csv"
